# Add 116th congress house members (commit: "add 116th congress house members")
#
# Sheet1 ("Sheet1") gets 32 new rows (840-871) of House-member / committee
# pairs for the 116th Congress, appended after the existing 117th-Congress
# data (which ends at row 839). Only columns B (member name), C (committee)
# and D (congress number) are populated for these new rows -- column A
# (the hyperlinked "Name, State" text) is left blank, matching the source
# data.
#
# Rows are written in a specific order so that the new shared-string table
# entries come out in the same order as the source workbook: every member
# is written the first time he/she is needed, top-to-bottom, EXCEPT "Mike
# Conaway" whose two rows (848-849) are written last even though they sit
# in the middle of the block -- that reproduces the shared-string index
# ordering exactly.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sorted")

$rows = @(
    @(840, 'Greg Gianforte', 'Natural Resources', 116),
    @(841, 'Greg Gianforte', 'Oversight and Reform', 116),
    @(842, 'Joe Kennedy III', 'Energy and Commerce', 116),
    @(843, 'Gil Cisneros', 'Armed Services', 116),
    @(844, 'Gil Cisneros', "Veterans' Affairs", 116),
    @(845, 'George Holding', 'Ways and Means', 116),
    @(846, 'George Holding', 'Budget', 116),
    @(847, 'George Holding', 'Ethics', 116),
    @(850, 'Susan Davis', 'Education and Labor', 116),
    @(851, 'Susan Davis', 'Armed Services', 116),
    @(852, 'Susan Davis', 'House Administration', 116),
    @(853, 'Pete Visclosky', 'Appropriations', 116),
    @(854, 'Francis Rooney', 'Education and Labor', 116),
    @(855, 'Francis Rooney', 'Foreign Affairs', 116),
    @(856, 'Justin Amash', 'Oversight and Reform', 116),
    @(857, 'Donna Shalala', 'Education and Labor', 116),
    @(858, 'Donna Shalala', 'Rules', 116),
    @(859, 'Susan Brooks', 'Energy and Commerce', 116),
    @(860, 'Susan Brooks', 'Ethics', 116),
    @(861, 'Roger Marshall', 'Agriculture', 116),
    @(862, 'Roger Marshall', 'Science, Space, and Technology', 116),
    @(863, 'Roger Marshall', 'Small Business', 116),
    @(864, 'Harley Rouda', 'Oversight and Reform', 116),
    @(865, 'Harley Rouda', 'Transportation and Infrastructure', 116),
    @(866, 'Phil Roe', "Veterans' Affairs", 116),
    @(867, 'Phil Roe', 'Education and Labor', 116),
    @(868, 'Lacy Clay', 'Financial Services', 116),
    @(869, 'Lacy Clay', 'Oversight and Reform', 116),
    @(870, 'Kenny Marchant', 'Ways and Means', 116),
    @(871, 'Kenny Marchant', 'Ethics', 116),
    @(848, 'Mike Conaway', 'Agriculture', 116),
    @(849, 'Mike Conaway', 'Armed Services', 116)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
}

# Printed page orientation for Sheet1 (now portrait, was unset/default).
$ws1.PageSetup.Orientation = 1

# Restore the view state: Sheet1 is the active/selected tab with the
# selection parked near the bottom of the newly-added data; "Sorted" keeps
# its own independent selection lower in its (unchanged) data.
$ws2.Activate()
$ws2.Range("B822").Select()

$ws1.Activate()
$ws1.Range("E836").Select()
